# Project_plan.xlsx update: mark "TÌNH TRẠNG" (F) column as "Hoàn thành"
# for every detail row (skips the section/group header rows), and move
# the active selection to G17 with the view scrolled near row 10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProjectPlan")

$detailRows = @(3,4,5,6,8,9,10,11,12,13,14,15,16,18,19,20,21,22,23,24,25)
foreach ($r in $detailRows) {
    $ws.Cells.Item($r, 6).Value = "Hoàn thành"
}

$ws.Activate()
$ws.Range("A10").Select()
$ws.Range("G17").Select()
